$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) hold text-formatted numbers (e.g. "606.88", "1.00").
# Pre-set NumberFormat to Text so Excel stores the literal string instead of
# silently coercing it to a numeric value (which would drop trailing zeros /
# introduce floating-point artifacts).
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D10", "D11", "D13", "D15", "D16", "D17", "D19", "D20", "D22", "D23", "D24", "D25", "D29", "D32", "D33", "D34", "D35", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $priceCells) {
  $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.298.70'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '3.567.43'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '606.88'
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").Value = '144.85'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").Value = '3.566.78'
$ws.Range("E7").Value = '  +0.71%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").Value = '  +2.02%  '
$ws.Range("D10").Value = '0.136'
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("D11").Value = '7.80'
$ws.Range("E11").Value = '  -2.66%  '
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("D13").Value = '4.175.00'
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").Value = '30.37'
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").Value = '3.561.70'
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("D17").Value = '66.371.01'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").Value = '11.51'
$ws.Range("E19").Value = '  +5.11%  '
$ws.Range("D20").Value = '6.23'
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("E21").Value = '  -1.42%  '
$ws.Range("D22").Value = '431.72'
$ws.Range("E22").Value = '  +1.27%  '
$ws.Range("D23").Value = '0.613'
$ws.Range("E23").Value = '  +1.71%  '
$ws.Range("D24").Value = '79.73'
$ws.Range("E24").Value = '  +1.17%  '
$ws.Range("D25").Value = '3.711.21'
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("E27").Value = '  -1.03%  '
$ws.Range("E28").Value = '  +1.08%  '
$ws.Range("D29").Value = '9.16'
$ws.Range("E29").Value = '  -2.02%  '
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("B32").Value = 'RenzoRestakedETH'
$ws.Range("C32").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D32").Value = '3.564.15'
$ws.Range("E32").Value = '  +0.89%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '25.47'
$ws.Range("E33").Value = '  +0.37%  '
$ws.Range("D34").Value = '1.46'
$ws.Range("E34").Value = '  -2.53%  '
$ws.Range("D35").Value = '0.152'
$ws.Range("E35").Value = '  -5.34%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  -2.11%  '
$ws.Range("D39").Value = '5.62'
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").Value = '175.55'
$ws.Range("E40").Value = '  +2.97%  '
$ws.Range("D41").Value = '0.0852'
$ws.Range("E41").Value = '  -1.37%  '
$ws.Range("D42").Value = '5.21'
$ws.Range("E42").Value = '  +0.60%  '
$ws.Range("D43").Value = '0.889'
$ws.Range("E43").Value = '  -0.53%  '
$ws.Range("E44").Value = '  +2.63%  '
$ws.Range("D45").Value = '45.99'
$ws.Range("E45").Value = '  +1.55%  '
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").Value = '2.52'
$ws.Range("E47").Value = '  +4.28%  '
$ws.Range("D48").Value = '1.20'
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("D49").Value = '25.20'
$ws.Range("E49").Value = '  -3.82%  '
$ws.Range("D50").Value = '7.14'
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").Value = '23.39'
$ws.Range("E51").Value = '  +3.37%  '
